$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4111.7646
$ws.Range("J17").Value = 4111.7646
$ws.Range("L17").Value = 12335.2938
$ws.Range("N17").Value = -12671.2938
$ws.Range("H64").Value = 8302.200000000001
$ws.Range("H67").Value = 8302.200000000001
$ws.Range("H80").Value = 1769.5
$ws.Range("I80").Value = 916.6
$ws.Range("J80").Value = 2157.182
$ws.Range("K80").Value = 2749.8
$ws.Range("L80").Value = 6471.545999999999
$ws.Range("M80").Value = -1751.8
$ws.Range("N80").Value = -8467.545999999998
$ws.Range("H83").Value = 1769.5
$ws.Range("I83").Value = 916.6
$ws.Range("J83").Value = 2157.182
$ws.Range("K83").Value = 8249.4
$ws.Range("L83").Value = 19414.638
$ws.Range("M83").Value = -3257.4
$ws.Range("N83").Value = -29398.638
$ws.Range("H132").Value = 1103.2413
$ws.Range("I132").Value = 874.7692
$ws.Range("J132").Value = 3083.3333
$ws.Range("K132").Value = 2624.3076
$ws.Range("L132").Value = 9249.999899999999
$ws.Range("M132").Value = -94.30760000000009
$ws.Range("N132").Value = -14309.9999
$ws.Range("H137").Value = 2565.5557
$ws.Range("I137").Value = 2084.6191
$ws.Range("J137").Value = 4248.8335
$ws.Range("K137").Value = 6253.8573
$ws.Range("L137").Value = 12746.5005
$ws.Range("M137").Value = -3703.8573
$ws.Range("N137").Value = -17846.5005

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2526176.5
$ws.Range("I2").Value = 3368053.5
$ws.Range("J2").Value = 545.44446
$ws.Range("K2").Value = 3368053.5
$ws.Range("L2").Value = 545.44446
$ws.Range("M2").Value = -3367940.5
$ws.Range("N2").Value = -771.44446
$ws.Range("H32").Value = 29405.658
$ws.Range("I32").Value = 29821.25
$ws.Range("K32").Value = 29821.25
$ws.Range("M32").Value = -29534.25
$ws.Range("H74").Value = 5094
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 5094
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 5094
$ws.Range("N74").Value = -6842
$ws.Range("H77").Value = 5094
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 5094
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 25470
$ws.Range("N77").Value = -34206
$ws.Range("H110").Value = 7354124
$ws.Range("I110").Value = 8929282
$ws.Range("K110").Value = 8929282
$ws.Range("M110").Value = -8927237
$ws.Range("H116").Value = 2526176.5
$ws.Range("I116").Value = 3368053.5
$ws.Range("J116").Value = 545.44446
$ws.Range("K116").Value = 3368053.5
$ws.Range("L116").Value = 545.44446
$ws.Range("M116").Value = -3365759.5
$ws.Range("N116").Value = -5133.44446
$ws.Range("M74").ClearContents()
$ws.Range("M77").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2526176.5
$ws.Range("I3").Value = 3368053.5
$ws.Range("J3").Value = 545.44446
$ws.Range("K3").Value = 3368053.5
$ws.Range("L3").Value = 545.44446
$ws.Range("M3").Value = -3367939.5
$ws.Range("N3").Value = -773.44446
$ws.Range("H22").Value = 2476.3333
$ws.Range("I22").Value = 1660.875
$ws.Range("K22").Value = 1660.875
$ws.Range("M22").Value = -1487.875
$ws.Range("H70").Value = 399999.94
$ws.Range("J70").Value = 399999.94
$ws.Range("L70").Value = 399999.94
$ws.Range("N70").Value = -400585.94
$ws.Range("H73").Value = 399999.94
$ws.Range("J73").Value = 399999.94
$ws.Range("L73").Value = 399999.94
$ws.Range("N73").Value = -402027.94
$ws.Range("H86").Value = 144552.14
$ws.Range("J86").Value = 251488.25
$ws.Range("L86").Value = 251488.25
$ws.Range("N86").Value = -253734.25
$ws.Range("H89").Value = 144552.14
$ws.Range("J89").Value = 251488.25
$ws.Range("L89").Value = 1257441.25
$ws.Range("N89").Value = -1268673.25
$ws.Range("H106").Value = 50000
$ws.Range("J106").Value = 50000
$ws.Range("L106").Value = 50000
$ws.Range("N106").Value = -52524

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 23953.652
$ws.Range("J9").Value = 23953.652
$ws.Range("L9").Value = 23953.652
$ws.Range("N9").Value = -24289.652
$ws.Range("H58").Value = 6098.0586
$ws.Range("I58").Value = 6133.1665
$ws.Range("K58").Value = 6133.1665
$ws.Range("M58").Value = -5930.1665
$ws.Range("H62").Value = 13724.25
$ws.Range("I62").Value = 4965.6665
$ws.Range("K62").Value = 4965.6665
$ws.Range("M62").Value = -4341.6665
$ws.Range("H65").Value = 13724.25
$ws.Range("I65").Value = 4965.6665
$ws.Range("K65").Value = 24828.3325
$ws.Range("M65").Value = -21708.3325
$ws.Range("H122").Value = 79621.16
$ws.Range("I122").Value = 92279.55
$ws.Range("K122").Value = 276838.65
$ws.Range("M122").Value = -274388.65
$ws.Range("H132").Value = 63695
$ws.Range("I132").Value = 7203.6
$ws.Range("K132").Value = 21610.8
$ws.Range("M132").Value = -19080.8
$ws.Range("H136").Value = 6098.0586
$ws.Range("I136").Value = 6133.1665
$ws.Range("K136").Value = 18399.4995
$ws.Range("M136").Value = -15849.4995

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 250775
$ws.Range("I47").Value = 250775
$ws.Range("K47").Value = 752325
$ws.Range("M47").Value = -751894
$ws.Range("H80").Value = 3375
$ws.Range("J80").Value = 2833.3333
$ws.Range("L80").Value = 8499.999899999999
$ws.Range("N80").Value = -10371.9999
$ws.Range("H83").Value = 3375
$ws.Range("J83").Value = 2833.3333
$ws.Range("L83").Value = 25499.9997
$ws.Range("N83").Value = -34859.9997
$ws.Range("H140").Value = 1173.0769
$ws.Range("J140").Value = 1241.8182
$ws.Range("L140").Value = 3725.4546
$ws.Range("N140").Value = -14085.4546

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 250500
$ws.Range("I12").Value = 500000
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 500000
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = -499860
$ws.Range("N12").Value = -1280
$ws.Range("H70").Value = 2953
$ws.Range("I70").Value = 2953
$ws.Range("K70").Value = 2953
$ws.Range("M70").Value = -2683
$ws.Range("H73").Value = 2953
$ws.Range("I73").Value = 2953
$ws.Range("K73").Value = 2953
$ws.Range("M73").Value = -2017
$ws.Range("H95").Value = 34248.5
$ws.Range("J95").Value = 34248.5
$ws.Range("L95").Value = 34248.5
$ws.Range("N95").Value = -39740.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 15693.625
$ws.Range("I40").Value = 14758.083
$ws.Range("J40").Value = 18500.25
$ws.Range("K40").Value = 14758.083
$ws.Range("L40").Value = 18500.25
$ws.Range("M40").Value = -14622.083
$ws.Range("N40").Value = -18772.25
$ws.Range("H61").Value = 4275.375
$ws.Range("I61").Value = 3799.75
$ws.Range("K61").Value = 3799.75
$ws.Range("M61").Value = -3597.75
$ws.Range("H68").Value = 2664.1428
$ws.Range("I68").Value = 2244.2222
$ws.Range("K68").Value = 2244.2222
$ws.Range("M68").Value = -1495.2222
$ws.Range("H71").Value = 2664.1428
$ws.Range("I71").Value = 2244.2222
$ws.Range("K71").Value = 11221.111
$ws.Range("M71").Value = -7477.111000000001
$ws.Range("H82").Value = 927.875
$ws.Range("I82").Value = 922.1111
$ws.Range("J82").Value = 935.2857
$ws.Range("K82").Value = 922.1111
$ws.Range("L82").Value = 935.2857
$ws.Range("M82").Value = -561.1111
$ws.Range("N82").Value = -1657.2857
$ws.Range("H85").Value = 927.875
$ws.Range("I85").Value = 922.1111
$ws.Range("J85").Value = 935.2857
$ws.Range("K85").Value = 922.1111
$ws.Range("L85").Value = 935.2857
$ws.Range("M85").Value = 325.8889
$ws.Range("N85").Value = -3431.2857
$ws.Range("H106").Value = 20954.857
$ws.Range("J106").Value = 20954.857
$ws.Range("L106").Value = 20954.857
$ws.Range("N106").Value = -23478.857
$ws.Range("H113").Value = 4275.375
$ws.Range("I113").Value = 3799.75
$ws.Range("K113").Value = 3799.75
$ws.Range("M113").Value = -1629.75
$ws.Range("H122").Value = 6250
$ws.Range("I122").Value = 6250
$ws.Range("K122").Value = 18750
$ws.Range("M122").Value = -16300
$ws.Range("H132").Value = 6368.25
$ws.Range("I132").Value = 2500
$ws.Range("J132").Value = 7657.6665
$ws.Range("K132").Value = 7500
$ws.Range("L132").Value = 22972.9995
$ws.Range("M132").Value = -4970
$ws.Range("N132").Value = -28032.9995
$ws.Range("H134").Value = 63080.832
$ws.Range("J134").Value = 63080.832
$ws.Range("L134").Value = 63080.832
$ws.Range("N134").Value = -73220.83199999999
$ws.Range("H136").Value = 4723.75
$ws.Range("I136").Value = 4660
$ws.Range("J136").Value = 4835.3125
$ws.Range("K136").Value = 13980
$ws.Range("L136").Value = 14505.9375
$ws.Range("M136").Value = -11430
$ws.Range("N136").Value = -19605.9375

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 75438.336
$ws.Range("J105").Value = 75438.336
$ws.Range("L105").Value = 75438.336
$ws.Range("N105").Value = -82426.336
$ws.Range("H122").Value = 2992.0476
$ws.Range("I122").Value = 2964.8948
$ws.Range("J122").Value = 3250
$ws.Range("K122").Value = 8894.6844
$ws.Range("L122").Value = 9750
$ws.Range("M122").Value = -6444.6844
$ws.Range("N122").Value = -14650
$ws.Range("H132").Value = 7420.7144
$ws.Range("I132").Value = 6717.273
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 20151.819
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -17621.819
$ws.Range("N132").Value = -35060
$ws.Range("H136").Value = 6082.6665
$ws.Range("J136").Value = 8807.714
$ws.Range("L136").Value = 26423.142
